$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Placeholder1 -> Placeholder
$ws.Range("A2").Value = "Stimuli/Placeholder.wav"
$ws.Range("B2").Value = "Stimuli/trigger_Placeholder.wav"

# Row 5: Harmonica3 -> Harmonica4
$ws.Range("A5").Value = "Stimuli/Harmonica4.wav"
$ws.Range("B5").Value = "Stimuli/trigger_Harmonica4.wav"

# Row 6: new Harmonica5 (was Placeholder2)
$ws.Range("A6").Value = "Stimuli/Harmonica5.wav"
$ws.Range("B6").Value = "Stimuli/trigger_Harmonica5.wav"

# Row 7: new Harmonica6 (was Placeholder3)
$ws.Range("A7").Value = "Stimuli/Harmonica6.wav"
$ws.Range("B7").Value = "Stimuli/trigger_Harmonica6.wav"

# Row 8: Placeholder2 shifted down (was Placeholder4)
$ws.Range("A8").Value = "Stimuli/Placeholder2.wav"
$ws.Range("B8").Value = "Stimuli/trigger_Placeholder2.wav"
